$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1778.3334
$ws.Range("I19").Value = 1102.875
$ws.Range("J19").Value = 2318.7
$ws.Range("K19").Value = 1102.875
$ws.Range("L19").Value = 2318.7
$ws.Range("M19").Value = -927.875
$ws.Range("N19").Value = -2668.7
$ws.Range("H43").Value = 7589.231
$ws.Range("J43").Value = 8765
$ws.Range("L43").Value = 8765
$ws.Range("N43").Value = -8903
$ws.Range("H141").Value = 2082.8215
$ws.Range("I141").Value = 2012.2693
$ws.Range("K141").Value = 6036.8079
$ws.Range("M141").Value = -856.8078999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2014.05
$ws.Range("I45").Value = 1384.1538
$ws.Range("J45").Value = 3183.8572
$ws.Range("K45").Value = 1384.1538
$ws.Range("L45").Value = 3183.8572
$ws.Range("M45").Value = -1007.1538
$ws.Range("N45").Value = -3937.8572
$ws.Range("H61").Value = 8794132
$ws.Range("I61").Value = 10207995
$ws.Range("J61").Value = 134221.5
$ws.Range("K61").Value = 10207995
$ws.Range("L61").Value = 134221.5
$ws.Range("M61").Value = -10207783
$ws.Range("N61").Value = -134645.5
$ws.Range("H88").Value = 2671
$ws.Range("J88").Value = 2507
$ws.Range("L88").Value = 2507
$ws.Range("N88").Value = -3319
$ws.Range("H91").Value = 2671
$ws.Range("J91").Value = 2507
$ws.Range("L91").Value = 2507
$ws.Range("N91").Value = -5315
$ws.Range("H119").Value = 734996.7
$ws.Range("J119").Value = 102495
$ws.Range("L119").Value = 102495
$ws.Range("N119").Value = -112171
$ws.Range("H122").Value = 1299.55
$ws.Range("I122").Value = 1117.2354
$ws.Range("K122").Value = 3351.7062
$ws.Range("M122").Value = -901.7062000000001
$ws.Range("H132").Value = 5875.143
$ws.Range("I132").Value = 2471.2
$ws.Range("K132").Value = 7413.599999999999
$ws.Range("M132").Value = -4883.599999999999
$ws.Range("H136").Value = 8794132
$ws.Range("I136").Value = 10207995
$ws.Range("J136").Value = 134221.5
$ws.Range("K136").Value = 30623985
$ws.Range("L136").Value = 402664.5
$ws.Range("M136").Value = -30621435
$ws.Range("N136").Value = -407764.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3466.1052
$ws.Range("I20").Value = 3079.9412
$ws.Range("K20").Value = 3079.9412
$ws.Range("M20").Value = -2832.9412
$ws.Range("H86").Value = 2796.3635
$ws.Range("I86").Value = 2574
$ws.Range("J86").Value = 3797
$ws.Range("K86").Value = 2574
$ws.Range("L86").Value = 3797
$ws.Range("M86").Value = -1451
$ws.Range("N86").Value = -6043
$ws.Range("H89").Value = 2796.3635
$ws.Range("I89").Value = 2574
$ws.Range("J89").Value = 3797
$ws.Range("K89").Value = 12870
$ws.Range("L89").Value = 18985
$ws.Range("M89").Value = -7254
$ws.Range("N89").Value = -30217
$ws.Range("H104").Value = 68000
$ws.Range("J104").Value = 68000
$ws.Range("L104").Value = 68000
$ws.Range("N104").Value = -74988
$ws.Range("H105").Value = 2742.625
$ws.Range("I105").Value = 2656.8333
$ws.Range("K105").Value = 2656.8333
$ws.Range("M105").Value = -909.8332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17833.334
$ws.Range("I41").Value = 9750
$ws.Range("J41").Value = 34000
$ws.Range("K41").Value = 9750
$ws.Range("L41").Value = 34000
$ws.Range("M41").Value = -9322
$ws.Range("N41").Value = -34856
$ws.Range("H59").Value = 69714.28999999999
$ws.Range("J59").Value = 69714.28999999999
$ws.Range("L59").Value = 69714.28999999999
$ws.Range("N59").Value = -72004.28999999999
$ws.Range("H60").Value = 4093
$ws.Range("I60").Value = 4093
$ws.Range("K60").Value = 4093
$ws.Range("M60").Value = -3582
$ws.Range("H86").Value = 6184.6665
$ws.Range("I86").Value = 6027.5
$ws.Range("J86").Value = 6499
$ws.Range("K86").Value = 6027.5
$ws.Range("L86").Value = 6499
$ws.Range("M86").Value = -4904.5
$ws.Range("N86").Value = -8745
$ws.Range("H89").Value = 6184.6665
$ws.Range("I89").Value = 6027.5
$ws.Range("J89").Value = 6499
$ws.Range("K89").Value = 30137.5
$ws.Range("L89").Value = 32495
$ws.Range("M89").Value = -24521.5
$ws.Range("N89").Value = -43727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 67.25
$ws.Range("J23").Value = 67.25
$ws.Range("L23").Value = 201.75
$ws.Range("N23").Value = -671.75
$ws.Range("H104").Value = 4874.75
$ws.Range("I104").Value = 4749.5
$ws.Range("K104").Value = 14248.5
$ws.Range("M104").Value = -11627.5
$ws.Range("H139").Value = 2684.5417
$ws.Range("I139").Value = 3632.25
$ws.Range("J139").Value = 2495
$ws.Range("K139").Value = 10896.75
$ws.Range("L139").Value = 7485
$ws.Range("M139").Value = -5756.75
$ws.Range("N139").Value = -17765

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5163
$ws.Range("I70").Value = 5163
$ws.Range("K70").Value = 5163
$ws.Range("M70").Value = -4893
$ws.Range("H73").Value = 5163
$ws.Range("I73").Value = 5163
$ws.Range("K73").Value = 5163
$ws.Range("M73").Value = -4227
$ws.Range("H92").Value = 31167.166
$ws.Range("J92").Value = 31167.166
$ws.Range("L92").Value = 31167.166
$ws.Range("N92").Value = -34911.166
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 51132.5
$ws.Range("J109").Value = 51132.5
$ws.Range("L109").Value = 51132.5
$ws.Range("N109").Value = -53212.5
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4148.355
$ws.Range("I40").Value = 3337.524
$ws.Range("J40").Value = 5851.1
$ws.Range("K40").Value = 3337.524
$ws.Range("L40").Value = 5851.1
$ws.Range("M40").Value = -3201.524
$ws.Range("N40").Value = -6123.1
$ws.Range("H43").Value = 2746630.8
$ws.Range("I43").Value = 4038200
$ws.Range("K43").Value = 4038200
$ws.Range("M43").Value = -4038007
$ws.Range("H55").Value = 43478700
$ws.Range("I55").Value = 62500440
$ws.Range("J55").Value = 440.42856
$ws.Range("K55").Value = 62500440
$ws.Range("L55").Value = 440.42856
$ws.Range("M55").Value = -62500267
$ws.Range("N55").Value = -786.4285600000001
$ws.Range("H61").Value = 1430
$ws.Range("I61").Value = 1430
$ws.Range("K61").Value = 1430
$ws.Range("M61").Value = -1228
$ws.Range("H82").Value = 859.1
$ws.Range("I82").Value = 599
$ws.Range("J82").Value = 1249.25
$ws.Range("K82").Value = 599
$ws.Range("L82").Value = 1249.25
$ws.Range("M82").Value = -238
$ws.Range("N82").Value = -1971.25
$ws.Range("H85").Value = 859.1
$ws.Range("I85").Value = 599
$ws.Range("J85").Value = 1249.25
$ws.Range("K85").Value = 599
$ws.Range("L85").Value = 1249.25
$ws.Range("M85").Value = 649
$ws.Range("N85").Value = -3745.25
$ws.Range("H113").Value = 1430
$ws.Range("I113").Value = 1430
$ws.Range("K113").Value = 1430
$ws.Range("M113").Value = 740
$ws.Range("H119").Value = 99995
$ws.Range("J119").Value = 99995
$ws.Range("L119").Value = 99995
$ws.Range("N119").Value = -109671
$ws.Range("H122").Value = 4802.3105
$ws.Range("I122").Value = 4569.0557
$ws.Range("K122").Value = 13707.1671
$ws.Range("M122").Value = -11257.1671
$ws.Range("H132").Value = 380919.7
$ws.Range("I132").Value = 358809.5
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 1076428.5
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -1073898.5
$ws.Range("N132").Value = -3005075
$ws.Range("H136").Value = 33697.902
$ws.Range("I136").Value = 5094.2666
$ws.Range("K136").Value = 15282.7998
$ws.Range("M136").Value = -12732.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 62450
$ws.Range("J64").Value = 62450
$ws.Range("L64").Value = 62450
$ws.Range("N64").Value = -62946
$ws.Range("H67").Value = 62450
$ws.Range("J67").Value = 62450
$ws.Range("L67").Value = 62450
$ws.Range("N67").Value = -64166
$ws.Range("H122").Value = 5765.75
$ws.Range("I122").Value = 3441
$ws.Range("K122").Value = 10323
$ws.Range("M122").Value = -7873
$ws.Range("H132").Value = 1565.6578
$ws.Range("I132").Value = 1352.6765
$ws.Range("J132").Value = 3376
$ws.Range("K132").Value = 4058.0295
$ws.Range("L132").Value = 10128
$ws.Range("M132").Value = -1528.0295
$ws.Range("N132").Value = -15188
$ws.Range("H136").Value = 696.125
$ws.Range("I136").Value = 652.7619
$ws.Range("K136").Value = 1958.2857
$ws.Range("M136").Value = 591.7143000000001
